# Update the "want to go" counts (column F) on the "展览" (rId1/sheet1)
# and "全部类型" (rId4/sheet4) worksheets, matching data refreshed from
# the source (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4790
    $ws.Range("F3").Value = 142
    $ws.Range("F4").Value = 833
}
